# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Every "Price" cell in column D is written with a leading apostrophe so Excel
# stores it as quote-prefixed TEXT (matching the sheet's existing text-based
# price column) instead of silently reinterpreting numeric-looking strings
# such as "598.46" as a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.253.06"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "'3.847.40"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'598.46"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "'167.52"
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("D7").Value = "'3.847.71"
$ws.Range("E7").Value = "  -2.12%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").Value = "'6.33"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "'0.462"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "'37.35"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "'4.490.83"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").Value = "'3.839.61"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").Value = "'68.403.65"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "'7.56"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "'18.19"
$ws.Range("E19").Value = "  +5.88%  "
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").Value = "'10.74"
$ws.Range("E21").Value = "  -4.07%  "
$ws.Range("D22").Value = "'472.40"
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").Value = "'0.735"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("E24").Value = "  -3.92%  "
$ws.Range("D25").Value = "'84.54"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("D27").Value = "'12.32"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("D28").Value = "'10.04"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").Value = "'3.996.36"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").Value = "'7.77"
$ws.Range("E32").Value = "  -1.32%  "
$ws.Range("D33").Value = "'2.31"
$ws.Range("E33").Value = "  -3.58%  "
$ws.Range("D34").Value = "'31.17"
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("D35").Value = "'3.819.48"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'5.95"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.30"
$ws.Range("E39").Value = "  +9.05%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "'1.01"
$ws.Range("E40").Value = "  -3.21%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").Value = "'0.315"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").Value = "'429.30"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("D44").Value = "'2.00"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'47.50"
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("D47").Value = "'8.61"
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("D48").Value = "'0.000275"
$ws.Range("E48").Value = "  +15.29%  "
$ws.Range("D49").Value = "'142.42"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "'0.0359"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "'39.15"
$ws.Range("E51").Value = "  -0.31%  "
